$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "bjj compression shorts",
    "compression running pants",
    "sit pad",
    "paintball pants padded",
    "knee pads breathable",
    "basketball compression tights for women",
    "elastic band black mountain",
    "lacrosse pads youth boys",
    "position pad",
    "knee length tights",
    "capri pants men",
    "knee pads volleyball black",
    "softball sliding pants women",
    "6ft basketball",
    "basketball shorts and pants",
    "work pants for men construction knee pads",
    "sliding shorts women softball",
    "fight shorts wrestling",
    "youth mesh leggings",
    "above the knee basketball shorts",
    "mens yoga leggings",
    "weight lifting pants for men",
    "cheap volleyball knee pads",
    "compression spandex",
    "yoga pants compression",
    "mens above the knee shorts",
    "mens running compression",
    "cycling pants mens",
    "knee sleeves basketball men",
    "softball gear for girls",
    "cold knee compression",
    "youth padded compression shorts",
    "yoga pants for men",
    "mens spandex tights",
    "softball protective gear",
    "soccer sliding shorts",
    "compression baseball shorts",
    "long shorts for men below knee",
    "padded leggings for cycling",
    "padded volleyball shorts",
    "hex squat",
    "youth padded sliding shorts",
    "knee sleeves bjj",
    "football pants pads adult",
    "work pants knee",
    "cold compression knee",
    "5 pad football girdle",
    "wrestling sleeve youth",
    "compression sports pants",
    "basketball tights for girls",
    "water pants",
    "spandex tights men",
    "boys compression pants black",
    "hockey tights",
    "youth hockey compression pants",
    "men leggings compression",
    "wrestling kneepads",
    "kneeling pad gym",
    "guard shorts",
    "padded compression shorts men",
    "softball pants youth",
    "spandex basketball shorts",
    "compression shorts men 5 pack",
    "shorts for men below knee",
    "mens gym leggings",
    "compression running leggings",
    "black mens basketball shorts",
    "knee pads impact",
    "paintball pads",
    "boys compression",
    "mens volleyball kneepads",
    "yoga knee pads",
    "knee work pads",
    "running capri",
    "paintball pants for men",
    "kneepad youth",
    "polyester capri pants",
    "man capri pants",
    "indoor baseball",
    "softball compression sleeve",
    "male pads",
    "high five girls softball pants",
    "outdoor hockey pants",
    "basketball floor",
    "basketball knee sleeve black",
    "long shorts for men below knee sports",
    "knee pads for adults",
    "hockey leggings",
    "volleyball long knee pads",
    "lacrosse shorts mens",
    "mens tights with pouch",
    "black short baseball pants",
    "lightweight knee pads",
    "mens compression pants cold",
    "knee shorts",
    "girls sliding pants",
    "knee pads for work for men",
    "youth padded leg sleeves for basketball",
    "gym knee compression",
    "compression football girdle",
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
